$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update forecast-error values: values shift one column to the left within each row
# (dropping the oldest evaluation), rows 2-6 gain a new value in column K,
# and rows 7-16 shrink by one column (their former last column is cleared).

# Row 2
$ws.Range("B2").Value = -0.5449296406117954
$ws.Range("C2").Value = -0.3099216406117954
$ws.Range("D2").Value = -0.6170576406117954
$ws.Range("E2").Value = -0.1063816406117954
$ws.Range("F2").Value = -0.3937486406117954
$ws.Range("G2").Value = -0.3657496406117954
$ws.Range("H2").Value = -0.2731016406117954
$ws.Range("I2").Value = -0.06012564061179543
$ws.Range("J2").Value = -0.5607856406117955
$ws.Range("K2").Value = -0.2710286406117954

# Row 3
$ws.Range("B3").Value = 0.2726728964218456
$ws.Range("C3").Value = -0.03446310357815441
$ws.Range("D3").Value = 0.4762128964218456
$ws.Range("E3").Value = 0.1888458964218456
$ws.Range("F3").Value = 0.2168448964218456
$ws.Range("G3").Value = 0.3094928964218456
$ws.Range("H3").Value = 0.5224688964218456
$ws.Range("I3").Value = 0.02180889642184558
$ws.Range("J3").Value = 0.3115658964218456
$ws.Range("K3").Value = 0.05367389642184558

# Row 4
$ws.Range("B4").Value = -0.1301381443875122
$ws.Range("C4").Value = 0.3805378556124878
$ws.Range("D4").Value = 0.09317085561248779
$ws.Range("E4").Value = 0.1211698556124878
$ws.Range("F4").Value = 0.2138178556124878
$ws.Range("G4").Value = 0.4267938556124878
$ws.Range("H4").Value = -0.07386614438751221
$ws.Range("I4").Value = 0.2158908556124878
$ws.Range("J4").Value = -0.04200114438751221
$ws.Range("K4").Value = 0.3403798556124878

# Row 5
$ws.Range("B5").Value = 0.6409568926112106
$ws.Range("C5").Value = 0.3535898926112106
$ws.Range("D5").Value = 0.3815888926112106
$ws.Range("E5").Value = 0.4742368926112106
$ws.Range("F5").Value = 0.6872128926112107
$ws.Range("G5").Value = 0.1865528926112106
$ws.Range("H5").Value = 0.4763098926112106
$ws.Range("I5").Value = 0.2184178926112106
$ws.Range("J5").Value = 0.6007988926112107
$ws.Range("K5").Value = 0.1895217986112106

# Row 6
$ws.Range("B6").Value = -0.4818651035472806
$ws.Range("C6").Value = -0.4538661035472806
$ws.Range("D6").Value = -0.3612181035472806
$ws.Range("E6").Value = -0.1482421035472806
$ws.Range("F6").Value = -0.6489021035472806
$ws.Range("G6").Value = -0.3591451035472806
$ws.Range("H6").Value = -0.6170371035472806
$ws.Range("I6").Value = -0.2346561035472806
$ws.Range("J6").Value = -0.6459331975472806
$ws.Range("K6").Value = -0.3589771035472806

# Row 7
$ws.Range("B7").Value = -0.3537867436446591
$ws.Range("C7").Value = -0.2611387436446591
$ws.Range("D7").Value = -0.04816274364465911
$ws.Range("E7").Value = -0.5488227436446591
$ws.Range("F7").Value = -0.2590657436446591
$ws.Range("G7").Value = -0.5169577436446591
$ws.Range("H7").Value = -0.1345767436446591
$ws.Range("I7").Value = -0.5458538376446591
$ws.Range("J7").Value = -0.2588977436446591

# Row 8
$ws.Range("B8").Value = 0.09264776243503714
$ws.Range("C8").Value = 0.3056237624350371
$ws.Range("D8").Value = -0.1950362375649629
$ws.Range("E8").Value = 0.09472076243503715
$ws.Range("F8").Value = -0.1631712375649629
$ws.Range("G8").Value = 0.2192097624350371
$ws.Range("H8").Value = -0.1920673315649629
$ws.Range("I8").Value = 0.09488876243503713

# Row 9
$ws.Range("B9").Value = 0.1574463720025918
$ws.Range("C9").Value = -0.3432136279974082
$ws.Range("D9").Value = -0.05345662799740819
$ws.Range("E9").Value = -0.3113486279974082
$ws.Range("F9").Value = 0.0710323720025918
$ws.Range("G9").Value = -0.3402447219974082
$ws.Range("H9").Value = -0.0532886279974082

# Row 10
$ws.Range("B10").Value = -0.5006596170015631
$ws.Range("C10").Value = -0.2109026170015632
$ws.Range("D10").Value = -0.4687946170015632
$ws.Range("E10").Value = -0.08641361700156319
$ws.Range("F10").Value = -0.4976907110015631
$ws.Range("G10").Value = -0.2107346170015632

# Row 11
$ws.Range("B11").Value = 0.2897568395245076
$ws.Range("C11").Value = 0.03186483952450761
$ws.Range("D11").Value = 0.4142458395245076
$ws.Range("E11").Value = 0.002968745524507627
$ws.Range("F11").Value = 0.2899248395245076

# Row 12
$ws.Range("B12").Value = -0.2484930410109615
$ws.Range("C12").Value = 0.1338879589890384
$ws.Range("D12").Value = -0.2773891350109615
$ws.Range("E12").Value = 0.009566958989038449

# Row 13
$ws.Range("B13").Value = 0.3058628168340501
$ws.Range("C13").Value = -0.1054142771659499
$ws.Range("D13").Value = 0.1815418168340501

# Row 14
$ws.Range("B14").Value = -0.4112768169122814
$ws.Range("C14").Value = -0.1243207229122814

# Row 15
$ws.Range("B15").Value = 0.4880092297750048

# Clear cells that no longer have data (row shrank / lost its only value)
$ws.Range("K7").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("G11").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("B16").ClearContents()
